$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 53. This shifts the existing rows 53-92
# down to 54-93 (carrying their values/formatting with them), and leaves
# a blank row 53 ready to receive the new week's data.
$ws.Rows.Item(53).Insert()

# Populate the newly inserted row 53 with the new record.
$ws.Cells.Item(53, 1).Value = 11
$ws.Cells.Item(53, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(53, 3).Value = "Bíobío"
$ws.Cells.Item(53, 4).Value = (Get-Date -Year 2022 -Month 1 -Day 5 -Hour 0 -Minute 0 -Second 0).Date
$ws.Cells.Item(53, 5).Value = 8
$ws.Cells.Item(53, 6).Value = "Fruta"
$ws.Cells.Item(53, 7).Value = 100103
$ws.Cells.Item(53, 8).Value = "Frutos de hueso (carozo)"
$ws.Cells.Item(53, 9).Value = 100103001
$ws.Cells.Item(53, 10).Value = "Cereza"
$ws.Cells.Item(53, 11).Value = "Lapins"
$ws.Cells.Item(53, 12).Value = "Primera"
$ws.Cells.Item(53, 13).Value = 250
$ws.Cells.Item(53, 14).Value = 5500
$ws.Cells.Item(53, 15).Value = 6000
$ws.Cells.Item(53, 16).Value = 5700
$ws.Cells.Item(53, 17).Value = "$/caja 10 kilos"
$ws.Cells.Item(53, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(53, 19).Value = 570
$ws.Cells.Item(53, 20).Value = 10
